# Input data read based on column headers
# Also added list of models and metric combinations that didn't converge on tab 2

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Sheet2 must be created from the ORIGINAL Sheet1 data (kVec/fVec/eVec/cVec
#    columns, B:E) before Sheet1 itself is edited.
# ---------------------------------------------------------------------------
$origB = @(2, 11, 2, 4, 3, 1, 1, 2, 4, 0, 4, 1, 3, 0)
$origC = @(1.3, 17.8, 5, 1.5, 1.5, 3, 3, 8, 30, 9, 25, 15, 15, 2)
$origD = @(0.05, 1, 0.19, 0.41, 0.32, 0.61, 0.32, 1.83, 3.01, 1.79, 3.17, 3.4, 4.2, 1.2)
$origE = @(0.5, 2.8, 1, 0.5, 0.5, 1, 0.5, 2.5, 3, 3, 6, 4, 4, 1)

# Move the selection on Sheet1 to its final resting place (B1) while it is
# still the active sheet, so the later Worksheets.Add() is what ends up
# owning the "active tab".
$ws1.Range("B1").Select()

# ---------------------------------------------------------------------------
# 2) Insert the new worksheet right after Sheet1 -> becomes "Sheet2" and is
#    activated automatically (matches activeTab="1" / tabSelected="1").
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)

$ws2.Cells.Item(1, 1).Value = "FC"
$ws2.Cells.Item(1, 2).Value = "fVec"
$ws2.Cells.Item(1, 3).Value = "eVec"
$ws2.Cells.Item(1, 4).Value = "cVec"

for ($i = 0; $i -lt 14; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 1).Value = $origB[$i]
    $ws2.Cells.Item($r, 2).Value = $origC[$i]
    $ws2.Cells.Item($r, 3).Value = $origD[$i]
    $ws2.Cells.Item($r, 4).Value = $origE[$i]
}

$ws2.Range("E1:E15").Select()

# ---------------------------------------------------------------------------
# 3) Sheet1: relabel the header row ("FC" takes B1, the rest of the headers
#    shift one slot over: fVec->C1, eVec->D1, cVec->E1, rVec->F1), keep
#    columns A-E's data as-is (column B keeps its original kVec numbers) and
#    refresh column F ("rVec") with its newly re-computed values.
# ---------------------------------------------------------------------------
$ws1.Cells.Item(1, 2).Value = "FC"
$ws1.Cells.Item(1, 3).Value = "fVec"
$ws1.Cells.Item(1, 4).Value = "eVec"
$ws1.Cells.Item(1, 5).Value = "cVec"
$ws1.Cells.Item(1, 6).Value = "rVec"

$newF_rVec = @(0.05, 0.34, 0.57999999999999996, 0.1, 0.03, 0.5, 0.23, 2, 0.05, 4, 3.6, 0.8, 4.3, 2)

for ($i = 0; $i -lt 14; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 6).Value = $newF_rVec[$i]
}

# NOTE: Sheet1's selection (B1) was already set above, before Sheet2 was
# added/activated. Re-selecting here would re-activate Sheet1 and steal the
# "active tab" status back from Sheet2, so we deliberately leave it alone.
